$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet: update Title and Date values
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(5, 2).Value = "NG-Imm Digital Signature"
$meta.Cells.Item(8, 2).Value = "2025-06-24T09:13:37+01:00"

# ---------------------------------------------------------------------
# 2. Elements sheet: content + structural changes
# ---------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# The root Extension element's "Short" text reused the same literal
# string as the document Title ("Digital Signature"), so it needs to
# track the rename too.
$elements.Cells.Item(2, 12).Value = "NG-Imm Digital Signature"

# -- Remove the autofilter (and its filter definitions) --------------
$elements.AutoFilterMode = $false

# -- Remove the conditional formatting rules (and backing dxfs) ------
$elements.Cells.FormatConditions.Delete()

# -- Un-hide the data rows (2-6) --------------------------------------
$elements.Rows.Item(2).Hidden = $false
$elements.Rows.Item(3).Hidden = $false
$elements.Rows.Item(4).Hidden = $false
$elements.Rows.Item(5).Hidden = $false
$elements.Rows.Item(6).Hidden = $false

# -- Update row 6 (Extension.value[x]) content so it absorbs the
#    information that used to live on the now-removed slice row (7) --
$elements.Cells.Item(6, 11).Value = "string`n"
$elements.Cells.Item(6, 12).Value = "Digital signature of the reporting officer"
$elements.Cells.Item(6, 13).Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
$elements.Cells.Item(6, 28).Value = ""
$elements.Cells.Item(6, 29).Value = ""
$elements.Cells.Item(6, 31).Value = ""
$elements.Cells.Item(6, 36).Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`n"

# -- Delete row 7 (the valueSignature slice row is merged into row 6) -
$elements.Rows.Item(7).Delete()

# -- Column width adjustments ------------------------------------------
$elements.Columns.Item(1).ColumnWidth = 15.666666666666666
$elements.Columns.Item(3).ColumnWidth = 9.0
$elements.Columns.Item(11).ColumnWidth = 7.5

# ---------------------------------------------------------------------
# 3. Workbook-level: drop the _FilterDatabase defined name
# ---------------------------------------------------------------------
for ($i = $wb.Names.Count; $i -ge 1; $i--) {
    $wb.Names.Item($i).Delete()
}
